$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -21.98550000000001
    21 = -19.95809999999997
    23 = -20.20979999999998
    25 = -21.71359999999999
    53 = -21.7937
    57 = -22.37120000000001
    59 = -22.189
    69 = -21.62749999999999
    79 = -20.41040000000002
    83 = -21.8727
    93 = -21.289
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
